$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 307; this pushes the existing rows 307-315 down to
# 308-316 (matching the rest of the diff, which is purely the old rows
# shifted down by one) and extends the used range to A1:R316.
$ws.Rows.Item(307).Insert()

# Populate the newly inserted row 307 with the new record's data.
$ws.Cells.Item(307, 1).Value2 = 6
$ws.Cells.Item(307, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(307, 3).Value2 = "Metropolitana"
$ws.Cells.Item(307, 4).Value2 = 45075
$ws.Cells.Item(307, 5).Value2 = 13
$ws.Cells.Item(307, 6).Value2 = 100112001
$ws.Cells.Item(307, 7).Value2 = "Berenjena"
$ws.Cells.Item(307, 8).Value2 = "Sin especificar"
$ws.Cells.Item(307, 9).Value2 = "Primera"
$ws.Cells.Item(307, 10).Value2 = 410
$ws.Cells.Item(307, 11).Value2 = 5000
$ws.Cells.Item(307, 12).Value2 = 6000
$ws.Cells.Item(307, 13).Value2 = 5439
$ws.Cells.Item(307, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(307, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(307, 16).Value2 = 109
$ws.Cells.Item(307, 17).Value2 = 50
$ws.Cells.Item(307, 18).Value2 = "Hortaliza"
